$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 25
$ws.Range("D6").Select() | Out-Null
